# Rename the embedded logo pictures' docPr/name attributes.
#   Pearson logo inline pictures (descr = "...PearsonLogo.png") -> name "image2.png"
#   BTec logo inline picture     (descr = "BTec_Logo-Orange")   -> name "image1.jpg"
#
# The pictures live in the document's headers/footers (not the main
# story), so walk Sections -> Headers/Footers -> InlineShapes. Identify
# each picture by its (stable) AlternativeText/description, since that's
# unchanged by this edit, rather than by its current Name.

$d = $word.ActiveDocument

for ($secIdx = 1; $secIdx -le $d.Sections.Count; $secIdx++) {
    $sec = $d.Sections.Item($secIdx)

    for ($hfIdx = 1; $hfIdx -le 3; $hfIdx++) {
        $hdr = $sec.Headers.Item($hfIdx)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                $descr = $shp.AlternativeText
                if ($descr -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                } elseif ($descr -like "*PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }

        $ftr = $sec.Footers.Item($hfIdx)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                $descr = $shp.AlternativeText
                if ($descr -eq "BTec_Logo-Orange") {
                    $shp.Name = "image1.jpg"
                } elseif ($descr -like "*PearsonLogo.png") {
                    $shp.Name = "image2.png"
                }
            }
        }
    }
}

Write-Host "Renamed logo inline shapes."
